$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '86.832.28'
$ws.Range('E2').Value = '  +9.24%  '
$ws.Range('D3').Value = '3.333.92'
$ws.Range('E3').Value = '  +5.27%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '219.60'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +6.69%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '639.12'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.324'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +20.86%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.622'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +5.24%  '
$ws.Range('D10').Value = '3.343.82'
$ws.Range('E10').Value = '  +5.64%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.606'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.36%  '
$ws.Range('E12').Value = '  +8.33%  '
$ws.Range('E13').Value = '  +1.74%  '
$ws.Range('D14').Value = '3.932.07'
$ws.Range('E14').Value = '  +4.93%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '34.46'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +9.36%  '
$ws.Range('E16').Value = '  +2.16%  '
$ws.Range('D17').Value = '86.912.89'
$ws.Range('E17').Value = '  +9.30%  '
$ws.Range('D18').Value = '3.317.39'
$ws.Range('E18').Value = '  +5.12%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.70'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.61%  '
$ws.Range('E20').Value = '  +8.89%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '447.19'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.21%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.19'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.53%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.27'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.23%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.44'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +11.59%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.40'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +15.47%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.15'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +11.56%  '
$ws.Range('D27').Value = '3.436.85'
$ws.Range('E27').Value = '  +3.18%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '78.60'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.37%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0000132'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +8.78%  '
$ws.Range('E30').Value = '  +0.18%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.182'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +49.20%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '607.04'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +10.34%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '9.28'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.88%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.10%  '
$ws.Range('E35').Value = '  +5.51%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.06'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.12%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.150'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.58%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '23.44'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.46%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.53'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +15.82%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.420'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.06%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.17%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.19'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +19.19%  '
$ws.Range('E43').Value = '  +2.69%  '
$ws.Range('E44').Value = '  +14.61%  '
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '156.66'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.53%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '189.94'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.18%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '45.70'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.49%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.37'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.17%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.791'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.15%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '26.50'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +6.61%  '
